# Quarterly indexing esoteric bug-fix operation
#
# Column A holds quarter-start dates (1-Jan, 1-Apr, 1-Jul, 1-Oct as Excel
# serials). The fix re-indexes every date to the 15th of the month that
# follows the quarter-start month (e.g. 1988-07-01 -> 1988-08-15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    $d = [datetime]::FromOADate($serial)

    $month = $d.Month + 1
    $year = $d.Year
    if ($month -gt 12) {
        $month = $month - 12
        $year = $year + 1
    }

    $newDate = Get-Date -Year $year -Month $month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $newDate.ToOADate()
}
